$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: add the date for A6 and the "2.5" value for C6
$ws.Range("A6").Value = "2021-03-31"
$ws.Range("C6").Value = "2.5"

# Update the current selection to C7, matching the edited workbook's view state
$ws.Range("C7").Select()
